$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H43").Value = 1167.3334
$ws.Range("J43").Value = 1502
$ws.Range("L43").Value = 1502
$ws.Range("N43").Value = -1640
$ws.Range("H53").Value = 1476.091
$ws.Range("I53").Value = 1474.2
$ws.Range("J53").Value = 1495
$ws.Range("K53").Value = 1474.2
$ws.Range("L53").Value = 1495
$ws.Range("M53").Value = -837.2
$ws.Range("N53").Value = -2769
$ws.Range("H80").Value = 1191.2858
$ws.Range("J80").Value = 1674.5
$ws.Range("L80").Value = 5023.5
$ws.Range("N80").Value = -7019.5
$ws.Range("H83").Value = 1191.2858
$ws.Range("J83").Value = 1674.5
$ws.Range("L83").Value = 15070.5
$ws.Range("N83").Value = -25054.5
$ws.Range("H88").Value = 5000
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 5000
$ws.Range("N88").Value = -5812
$ws.Range("H91").Value = 5000
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 5000
$ws.Range("N91").Value = -7808
$ws.Range("H100").Value = 3000
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -4082
$ws.Range("H132").Value = 11911.25
$ws.Range("I132").Value = 12687.143
$ws.Range("K132").Value = 38061.429
$ws.Range("M132").Value = -35531.429
$ws.Range("H135").Value = 3235.5
$ws.Range("I135").Value = 3235.5
$ws.Range("K135").Value = 29119.5
$ws.Range("M135").Value = -26584.5
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 185.45454
$ws.Range("I5").Value = 177.88889
$ws.Range("K5").Value = 177.88889
$ws.Range("M5").Value = -65.88889
$ws.Range("H32").Value = 5853.4165
$ws.Range("I32").Value = 5853.4165
$ws.Range("K32").Value = 5853.4165
$ws.Range("M32").Value = -5566.4165
$ws.Range("H45").Value = 4400
$ws.Range("I45").Value = 4400
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 4400
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -4023
$ws.Range("N45").ClearContents()
$ws.Range("H74").Value = 2353.55
$ws.Range("I74").Value = 1385.2142
$ws.Range("K74").Value = 1385.2142
$ws.Range("M74").Value = -511.2141999999999
$ws.Range("H77").Value = 2353.55
$ws.Range("I77").Value = 1385.2142
$ws.Range("K77").Value = 6926.071
$ws.Range("M77").Value = -2558.071
$ws.Range("H101").Value = 45000.5
$ws.Range("J101").Value = 45000.5
$ws.Range("L101").Value = 45000.5
$ws.Range("N101").Value = -51490.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 185.45454
$ws.Range("I4").Value = 177.88889
$ws.Range("K4").Value = 177.88889
$ws.Range("M4").Value = -62.88889
$ws.Range("H105").Value = 41418.6
$ws.Range("I105").Value = 1773.75
$ws.Range("J105").Value = 199998
$ws.Range("K105").Value = 1773.75
$ws.Range("L105").Value = 199998
$ws.Range("M105").Value = -26.75
$ws.Range("N105").Value = -203492
$ws.Range("H134").Value = 3158.25
$ws.Range("I134").Value = 3227.182
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 9681.545999999998
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -7146.545999999998
$ws.Range("N134").Value = -12270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 226
$ws.Range("I10").Value = 164.22223
$ws.Range("J10").Value = 504
$ws.Range("K10").Value = 164.22223
$ws.Range("L10").Value = 504
$ws.Range("M10").Value = -25.22223
$ws.Range("N10").Value = -782
$ws.Range("H13").Value = 490
$ws.Range("J13").Value = 490
$ws.Range("L13").Value = 490
$ws.Range("N13").Value = -768
$ws.Range("H14").Value = 3331.1
$ws.Range("I14").Value = 3250
$ws.Range("J14").Value = 3452.75
$ws.Range("K14").Value = 3250
$ws.Range("L14").Value = 3452.75
$ws.Range("M14").Value = -3080
$ws.Range("N14").Value = -3792.75
$ws.Range("H28").Value = 26756.857
$ws.Range("J28").Value = 26756.857
$ws.Range("L28").Value = 26756.857
$ws.Range("N28").Value = -27246.857
$ws.Range("H31").Value = 3614.111
$ws.Range("I31").Value = 2129
$ws.Range("J31").Value = 4356.6665
$ws.Range("K31").Value = 2129
$ws.Range("L31").Value = 4356.6665
$ws.Range("M31").Value = -1834
$ws.Range("N31").Value = -4946.6665
$ws.Range("H34").Value = 3614.111
$ws.Range("I34").Value = 2129
$ws.Range("J34").Value = 4356.6665
$ws.Range("K34").Value = 2129
$ws.Range("L34").Value = 4356.6665
$ws.Range("M34").Value = -1927
$ws.Range("N34").Value = -4760.6665
$ws.Range("H74").Value = 56999.5
$ws.Range("J74").Value = 56999.5
$ws.Range("L74").Value = 56999.5
$ws.Range("N74").Value = -58747.5
$ws.Range("H77").Value = 56999.5
$ws.Range("J77").Value = 56999.5
$ws.Range("L77").Value = 170998.5
$ws.Range("N77").Value = -179734.5
$ws.Range("H132").Value = 3914.5
$ws.Range("J132").Value = 3244
$ws.Range("L132").Value = 9732
$ws.Range("N132").Value = -14792
$ws.Range("H134").Value = 712.5
$ws.Range("I134").Value = 712.5
$ws.Range("K134").Value = 2137.5
$ws.Range("M134").Value = 397.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1427.0667
$ws.Range("I4").Value = 773.4286
$ws.Range("K4").Value = 2320.2858
$ws.Range("M4").Value = -2208.2858
$ws.Range("H26").Value = 83
$ws.Range("I26").Value = 50
$ws.Range("J26").Value = 149
$ws.Range("K26").Value = 150
$ws.Range("L26").Value = 447
$ws.Range("M26").Value = 138
$ws.Range("N26").Value = -1023
$ws.Range("H60").Value = 430.81818
$ws.Range("J60").Value = 525
$ws.Range("L60").Value = 1575
$ws.Range("N60").Value = -2077
$ws.Range("H64").Value = 1012
$ws.Range("I64").Value = 1012
$ws.Range("K64").Value = 3036
$ws.Range("M64").Value = -2766
$ws.Range("H67").Value = 1012
$ws.Range("I67").Value = 1012
$ws.Range("K67").Value = 3036
$ws.Range("M67").Value = -2100
$ws.Range("H68").Value = 598.5625
$ws.Range("I68").Value = 636.6
$ws.Range("K68").Value = 1909.8
$ws.Range("M68").Value = -1098.8
$ws.Range("H69").Value = 2957
$ws.Range("J69").Value = 2957
$ws.Range("L69").Value = 8871
$ws.Range("N69").Value = -10493
$ws.Range("H71").Value = 598.5625
$ws.Range("I71").Value = 636.6
$ws.Range("K71").Value = 5729.400000000001
$ws.Range("M71").Value = -1673.400000000001
$ws.Range("H72").Value = 2957
$ws.Range("J72").Value = 2957
$ws.Range("L72").Value = 26613
$ws.Range("N72").Value = -34725

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 94415.27
$ws.Range("J14").Value = 35496.125
$ws.Range("L14").Value = 35496.125
$ws.Range("N14").Value = -35832.125
$ws.Range("H35").Value = 6000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 6000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 6000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -6596

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H40").Value = 2836.3333
$ws.Range("I40").Value = 2752
$ws.Range("J40").Value = 3005
$ws.Range("K40").Value = 2752
$ws.Range("L40").Value = 3005
$ws.Range("M40").Value = -2616
$ws.Range("N40").Value = -3277
$ws.Range("H82").Value = 1500.3334
$ws.Range("I82").Value = 1250.5
$ws.Range("J82").Value = 2000
$ws.Range("K82").Value = 1250.5
$ws.Range("L82").Value = 2000
$ws.Range("M82").Value = -889.5
$ws.Range("N82").Value = -2722
$ws.Range("H85").Value = 1500.3334
$ws.Range("I85").Value = 1250.5
$ws.Range("J85").Value = 2000
$ws.Range("K85").Value = 1250.5
$ws.Range("L85").Value = 2000
$ws.Range("M85").Value = -2.5
$ws.Range("N85").Value = -4496
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 1914.5
$ws.Range("I132").Value = 2004
$ws.Range("K132").Value = 6012
$ws.Range("M132").Value = -3482
$ws.Range("H136").Value = 1802.6666
$ws.Range("I136").Value = 1802.6666
$ws.Range("K136").Value = 5407.9998
$ws.Range("M136").Value = -2857.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 685
$ws.Range("I7").Value = 1025
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 1025
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = -912
$ws.Range("N7").Value = -231
$ws.Range("H113").Value = 4176.7144
$ws.Range("I113").Value = 1790.2
$ws.Range("J113").Value = 4922.5
$ws.Range("K113").Value = 5370.6
$ws.Range("L113").Value = 14767.5
$ws.Range("M113").Value = -3200.6
$ws.Range("N113").Value = -19107.5
$ws.Range("H122").Value = 855.7143
$ws.Range("I122").Value = 798.5
$ws.Range("K122").Value = 2395.5
$ws.Range("M122").Value = 54.5
